
# "Added deleting variables in screen 2"
#
# The Translation sheet (2nd sheet) held, among others, these rows (24-35):
#   24 <value>                                   (weight input value)
#   25 0                                         (weight default)
#   26 <value>                                   (height input value)
#   27 0                                         (height default)
#   28 Waga                                      (label "Weight")
#   29 Wzrost                                    (label "Height")
#   30 Pojemnosc\n butelki                       (label "Bottle capacity")
#   31 kg                                        (unit)
#   32 <value> l                                 (capacity value + unit)
#   33 Twoje zapotrzebowanie na wodę wynosi:     (result text)
#   34 cm                                        (unit)
#   35 l                                         (unit)
#
# The weight/height variables were removed from the screen, so rows 24, 25,
# 29 and 34 (weight value, weight default, "Wzrost" label and "cm" unit) are
# deleted. The remaining rows shift up. The capacity is now expressed in
# millilitres instead of litres ("<value> l" -> "<value> ml", "l" -> "ml"),
# and a new "x <value>" text is appended as a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Delete rows bottom-to-top so earlier row numbers stay valid.
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# After the deletions the surviving rows now occupy 24-31 in this order:
#   24 <value>                                   (was 26)
#   25 0                                         (was 27)
#   26 Waga                                      (was 28)
#   27 Pojemnosc\n butelki                       (was 30)
#   28 kg                                        (was 31)
#   29 <value> l  -> needs to become <value> ml  (was 32)
#   30 Twoje zapotrzebowanie na wodę wynosi:     (was 33)
#   31 l  -> needs to become ml                  (was 35)

$ws.Range("F29").Value = "<value> ml"
$ws.Range("F31").Value = "ml"

# Add the new row 32 with the new "x <value>" text, matching the style of
# the other rows in the table (same typography/alignment/direction columns).
$ws.Range("B32").Value = "SingleUseId43"
$ws.Range("C32").Value = "Default"
$ws.Range("D32").Value = "Left"
$ws.Range("E32").Value = "LTR"
$ws.Range("F32").Value = "x <value>"
